$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$br = [char]11

$t.Cell(1,1).Range.Text = "17 x 31" + $br + "  3    1" + $br + "  ----" + $br + "1|    |" + $br + "7|    |"
$t.Cell(1,2).Range.Text = "16 x 45" + $br + "  4    5" + $br + "  ----" + $br + "1|    |" + $br + "6|    |"
$t.Cell(1,3).Range.Text = "81 x 91" + $br + "  9    1" + $br + "  ----" + $br + "8|    |" + $br + "1|    |"
$t.Cell(2,1).Range.Text = "55 x 28" + $br + "  2    8" + $br + "  ----" + $br + "5|    |" + $br + "5|    |"
$t.Cell(2,2).Range.Text = "65 x 30" + $br + "  3    0" + $br + "  ----" + $br + "6|    |" + $br + "5|    |"
$t.Cell(2,3).Range.Text = "59 x 45" + $br + "  4    5" + $br + "  ----" + $br + "5|    |" + $br + "9|    |"
$t.Cell(3,1).Range.Text = "67 x 66" + $br + "  6    6" + $br + "  ----" + $br + "6|    |" + $br + "7|    |"
$t.Cell(3,2).Range.Text = "31 x 61" + $br + "  6    1" + $br + "  ----" + $br + "3|    |" + $br + "1|    |"
$t.Cell(3,3).Range.Text = "67 x 45" + $br + "  4    5" + $br + "  ----" + $br + "6|    |" + $br + "7|    |"
$t.Cell(4,1).Range.Text = "60 x 61" + $br + "  6    1" + $br + "  ----" + $br + "6|    |" + $br + "0|    |"
$t.Cell(4,2).Range.Text = "54 x 71" + $br + "  7    1" + $br + "  ----" + $br + "5|    |" + $br + "4|    |"
$t.Cell(4,3).Range.Text = "22 x 86" + $br + "  8    6" + $br + "  ----" + $br + "2|    |" + $br + "2|    |"
$t.Cell(5,1).Range.Text = "52 x 73" + $br + "  7    3" + $br + "  ----" + $br + "5|    |" + $br + "2|    |"
$t.Cell(5,2).Range.Text = "97 x 28" + $br + "  2    8" + $br + "  ----" + $br + "9|    |" + $br + "7|    |"
$t.Cell(5,3).Range.Text = "87 x 69" + $br + "  6    9" + $br + "  ----" + $br + "8|    |" + $br + "7|    |"
